# Customer workbook edit:
#  - "customer" sheet header E1: "address" -> "addressList"
#  - "customer" sheet: remove the two trailing blank formatting rows (9 & 10)
#  - Leave "customer" as the active sheet with E8 selected (matches author's
#    final cursor position after editing the header & trimming the rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer")

# Rename header from "address" to "addressList"
$ws.Range("E1").Value = "addressList"

# Remove the two empty trailing rows (9 and 10) that carried no data
$ws.Range("9:10").Delete()

# Make "customer" the active sheet/tab with E8 selected
$ws.Activate()
[void]$ws.Range("E8").Select()
